$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 397.66666
$ws.Range("J42").Value = 900
$ws.Range("L42").Value = 2700
$ws.Range("N42").Value = -3160
$ws.Range("H88").Value = 3426.8572
$ws.Range("J88").Value = 2498.1667
$ws.Range("L88").Value = 2498.1667
$ws.Range("N88").Value = -3310.1667
$ws.Range("H91").Value = 3426.8572
$ws.Range("J91").Value = 2498.1667
$ws.Range("L91").Value = 2498.1667
$ws.Range("N91").Value = -5306.1667
$ws.Range("H100").Value = 1799.8
$ws.Range("H106").Value = 2775.3845
$ws.Range("I106").Value = 3466.7778
$ws.Range("K106").Value = 3466.7778
$ws.Range("M106").Value = -2835.7778
$ws.Range("H116").Value = 9665.923000000001
$ws.Range("I116").Value = 14377.75
$ws.Range("J116").Value = 2127
$ws.Range("K116").Value = 14377.75
$ws.Range("L116").Value = 2127
$ws.Range("M116").Value = -10935.75
$ws.Range("N116").Value = -9011
$ws.Range("H132").Value = 844.6222
$ws.Range("I132").Value = 785.12195
$ws.Range("K132").Value = 2355.36585
$ws.Range("M132").Value = 174.6341499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4294.411
$ws.Range("I32").Value = 3328.403
$ws.Range("K32").Value = 3328.403
$ws.Range("M32").Value = -3041.403
$ws.Range("H61").Value = 6337.478
$ws.Range("I61").Value = 7297.6665
$ws.Range("J61").Value = 4537.125
$ws.Range("K61").Value = 7297.6665
$ws.Range("L61").Value = 4537.125
$ws.Range("M61").Value = -7085.6665
$ws.Range("N61").Value = -4961.125
$ws.Range("H74").Value = 767.8333
$ws.Range("I74").Value = 545.0303
$ws.Range("J74").Value = 3218.6667
$ws.Range("K74").Value = 545.0303
$ws.Range("L74").Value = 3218.6667
$ws.Range("M74").Value = 328.9697
$ws.Range("N74").Value = -4966.6667
$ws.Range("H77").Value = 767.8333
$ws.Range("I77").Value = 545.0303
$ws.Range("J77").Value = 3218.6667
$ws.Range("K77").Value = 2725.1515
$ws.Range("L77").Value = 16093.3335
$ws.Range("M77").Value = 1642.8485
$ws.Range("N77").Value = -24829.3335
$ws.Range("H97").Value = 1521.1538
$ws.Range("I97").Value = 1564.1111
$ws.Range("K97").Value = 1564.1111
$ws.Range("M97").Value = -1068.1111
$ws.Range("H102").Value = 1267.25
$ws.Range("I102").Value = 1095.3334
$ws.Range("J102").Value = 1439.1666
$ws.Range("K102").Value = 1095.3334
$ws.Range("L102").Value = 1439.1666
$ws.Range("M102").Value = 526.6666
$ws.Range("N102").Value = -4683.1666
$ws.Range("H132").Value = 1537.091
$ws.Range("I132").Value = 1184.6522
$ws.Range("K132").Value = 3553.9566
$ws.Range("M132").Value = -1023.9566
$ws.Range("H136").Value = 6337.478
$ws.Range("I136").Value = 7297.6665
$ws.Range("J136").Value = 4537.125
$ws.Range("K136").Value = 21892.9995
$ws.Range("L136").Value = 13611.375
$ws.Range("M136").Value = -19342.9995
$ws.Range("N136").Value = -18711.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2013.4286
$ws.Range("I94").Value = 819
$ws.Range("K94").Value = 819
$ws.Range("M94").Value = -368
$ws.Range("H105").Value = 2496.9048
$ws.Range("I105").Value = 2302
$ws.Range("K105").Value = 2302
$ws.Range("M105").Value = -555
$ws.Range("H134").Value = 7442.05
$ws.Range("I134").Value = 8121.0557
$ws.Range("K134").Value = 24363.1671
$ws.Range("M134").Value = -21828.1671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1959.0605
$ws.Range("I134").Value = 1748.3667
$ws.Range("K134").Value = 5245.1001
$ws.Range("M134").Value = -2710.1001
$ws.Range("H141").Value = 57613.43
$ws.Range("J141").Value = 55215.668
$ws.Range("L141").Value = 55215.668
$ws.Range("N141").Value = -65575.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -8888
$ws.Range("N3").ClearContents()
$ws.Range("H131").Value = 15319.667
$ws.Range("J131").Value = 15623.702
$ws.Range("L131").Value = 46871.106
$ws.Range("N131").Value = -56951.106
$ws.Range("H133").Value = 4052.5
$ws.Range("I133").Value = 2010
$ws.Range("J133").Value = 4733.3335
$ws.Range("K133").Value = 6030
$ws.Range("L133").Value = 14200.0005
$ws.Range("M133").Value = -970
$ws.Range("N133").Value = -24320.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 956.4643
$ws.Range("I97").Value = 1023.1053
$ws.Range("K97").Value = 1023.1053
$ws.Range("M97").Value = -527.1053000000001
$ws.Range("H132").Value = 1604937
$ws.Range("I132").Value = 2138371.5
$ws.Range("K132").Value = 6415114.5
$ws.Range("M132").Value = -6412584.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 65046
$ws.Range("J48").Value = 65046
$ws.Range("L48").Value = 65046
$ws.Range("N48").Value = -66368
$ws.Range("H61").Value = 3566.25
$ws.Range("I61").Value = 3293.8462
$ws.Range("J61").Value = 4072.1428
$ws.Range("K61").Value = 3293.8462
$ws.Range("L61").Value = 4072.1428
$ws.Range("M61").Value = -3091.8462
$ws.Range("N61").Value = -4476.1428
$ws.Range("H93").Value = 27778486
$ws.Range("I93").Value = 703
$ws.Range("K93").Value = 703
$ws.Range("M93").Value = 545
$ws.Range("H100").Value = 1297
$ws.Range("I100").Value = 1362.6666
$ws.Range("J100").Value = 1100
$ws.Range("K100").Value = 1362.6666
$ws.Range("L100").Value = 1100
$ws.Range("M100").Value = -821.6666
$ws.Range("N100").Value = -2182
$ws.Range("H113").Value = 3566.25
$ws.Range("I113").Value = 3293.8462
$ws.Range("J113").Value = 4072.1428
$ws.Range("K113").Value = 3293.8462
$ws.Range("L113").Value = 4072.1428
$ws.Range("M113").Value = -1123.8462
$ws.Range("N113").Value = -8412.1428
$ws.Range("H132").Value = 2778.0293
$ws.Range("I132").Value = 1626.25
$ws.Range("K132").Value = 4878.75
$ws.Range("M132").Value = -2348.75
$ws.Range("H136").Value = 1546.2162
$ws.Range("I136").Value = 1210.6666
$ws.Range("J136").Value = 2984.2856
$ws.Range("K136").Value = 3631.9998
$ws.Range("L136").Value = 8952.856800000001
$ws.Range("M136").Value = -1081.9998
$ws.Range("N136").Value = -14052.8568
